# Weekly update: insert a new data row at row 19 (pushing all existing
# price-history rows down by one) to record the latest "Puerro" price
# observation for Vega Central Mapocho de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 19:73 down to 20:74 and open up a fresh row 19.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row with the latest week's record.
$ws.Range("A19").Value = 9
$ws.Range("B19").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C19").Value = "Metropolitana"
$ws.Range("D19").Value = 44497
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = 100112005
$ws.Range("G19").Value = "Puerro"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 180
$ws.Range("K19").Value = 6000
$ws.Range("L19").Value = 7000
$ws.Range("M19").Value = 6556
$ws.Range("N19").Value = '$/paquete 20 unidades'
$ws.Range("O19").Value = "Provincia de Chacabuco"
$ws.Range("P19").Value = 328
$ws.Range("Q19").Value = 20
$ws.Range("R19").Value = "Hortaliza"
